$wb = $excel.ActiveWorkbook

# --- Add the new training-session row to the "s304" sheet ---
$ws304 = $wb.Worksheets.Item("s304")

# A3: date (copy A2's date-format style, then set the value)
$ws304.Range("A2").Copy()
$ws304.Range("A3").PasteSpecial(-4122)
$ws304.Range("A3").Value = "7/16/2015"

# C3:E3 plain text values (no special style)
$ws304.Range("C3").Value = "s304"
$ws304.Range("D3").Value = "Training"
$ws304.Range("E3").Value = "1 hour"

# F3: cost (copy F2's currency style, then set the value)
$ws304.Range("F2").Copy()
$ws304.Range("F3").PasteSpecial(-4122)
$ws304.Range("F3").Value = 10

# G3: paid? (copy G2's bold style, then set the value)
$ws304.Range("G2").Copy()
$ws304.Range("G3").PasteSpecial(-4122)
$ws304.Range("G3").Value = "No"

# H3: notes
$ws304.Range("H3").Value = "Added catch trial runs"

$ws304.Range("H43").Select()

# --- Add a new, empty worksheet named "s305" at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "s305"

# Match the page-margin conventions used by the other sheets in this workbook
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Keep "s304" as the active sheet/tab, as in the target workbook
$ws304.Activate()
